# Week 1 predictions workbook update:
#  - Highlight the cells whose model prediction matched the actual game
#    winner with a yellow fill (keeping existing borders).
#  - Add a new row 18 with each model's hit-rate (win %) formatted as a
#    percentage.
#  - Update the view (zoom level + active selection) to match the state
#    the workbook was left in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose prediction matched the actual winner -> mark with yellow fill.
# (Existing cell borders are preserved automatically since Interior.Color
# only touches the fill, not the border.)
$yellowCells = @(
    "C3", "D3", "E3", "F3", "G3",
    "C6", "E6", "F6",
    "C7", "D7", "E7", "F7",
    "D9",
    "C13", "D13", "E13", "F13",
    "C14", "D14", "E14", "F14", "G14",
    "C15", "F15",
    "G16",
    "G17"
)
foreach ($addr in $yellowCells) {
    $ws.Range($addr).Interior.Color = 65535
}

# New row 18: per-model win percentage across the week's games.
$ws.Range("C18").Value = 0.625
$ws.Range("D18").Value = 0.6875
$ws.Range("E18").Value = 0.6875
$ws.Range("F18").Value = 0.625
$ws.Range("G18").Value = 0.75

# Apply number formats in this order so the underlying style table matches:
# "0%" (Greg's Picks) is registered before "0.00%" (the model columns).
$ws.Range("G18").NumberFormat = "0%"
$ws.Range("C18:F18").NumberFormat = "0.00%"

# Restore the view state: zoomed out a bit and selection moved to F26.
$excel.ActiveWindow.Zoom = 115
[void]$ws.Range("F26").Select()
